$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("0.case")
$ws.Activate()

$ws.Range("A11").Value = "NSGA2"

$ws.Range("A13").Value = "ang_co"
$ws.Range("B13").Value = "X1"
$ws.Range("C13").Value = 15
$ws.Range("D13").Value = 25

$ws.Range("A14").Value = "deg_co"
$ws.Range("B14").Value = "X2"
$ws.Range("C14").Value = 90
$ws.Range("D14").Value = 150

$ws.Range("A15").Value = "bd"
$ws.Range("B15").Value = "X3"
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 4

$ws.Range("A16").Value = "bw"
$ws.Range("B16").Value = "X4"
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 1

$ws.Range("A17").Value = "bh"
$ws.Range("B17").Value = "X5"
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 4

$ws.Range("A18").Value = "bg"
$ws.Range("B18").Value = "X6"
$ws.Range("C18").Value = 1.5
$ws.Range("D18").Value = 2.5

$ws.Range("D20").Select()
